# EV Charger Cost.xlsx — record the state (Oregon) and the date this
# workbook revision was produced on the "About" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# New label next to the title, identifying the state this input file is for.
$ws.Range("B1").Value = "Oregon"

# Stamp the revision date (2022-03-11) with a standard short date format.
$ws.Range("C1").Value = 44631
$ws.Range("C1").NumberFormat = "mm-dd-yy"
